$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "307.44"
$ws.Range("E2").Value = "-0.10%"
$rng.Style = "Normal"

# Row 3
$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "39.84"
$ws.Range("E3").Value = "0.98%"
$rng.Style = "Normal"

# Row 4
$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.148"
$ws.Range("E4").Value = "0.70%"
$rng.Style = "Normal"

# Row 5
$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.08096"
$ws.Range("E5").Value = "-0.67%"
$rng.Style = "Normal"

# Row 6
$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "1.935"
$ws.Range("E6").Value = "-2.24%"
$rng.Style = "Normal"

# Row 7
$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "8.151"
$ws.Range("E7").Value = "3.19%"
$rng.Style = "Normal"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.9286"
$ws.Range("E8").Value = "0.02%"
$rng.Style = "Normal"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "0.1412"
$ws.Range("E9").Value = "0.26%"
$rng.Style = "Normal"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.1918"
$ws.Range("E10").Value = "-1.42%"
$rng.Style = "Normal"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.09077"
$ws.Range("E11").Value = "-1.94%"
$rng.Style = "Normal"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.03500"
$ws.Range("E12").Value = "-0.11%"
$rng.Style = "Normal"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.09816"
$ws.Range("E13").Value = "-0.12%"
$rng.Style = "Normal"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001391"
$ws.Range("E14").Value = "-0.94%"
$rng.Style = "Normal"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.005911"
$ws.Range("E15").Value = "0.17%"
$rng.Style = "Normal"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "3.947"
$ws.Range("E16").Value = "0.02%"
$rng.Style = "Normal"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "4.227"
$ws.Range("E17").Value = "1.14%"
$rng.Style = "Normal"

# Row 18
$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = "-2.06%"
$rng.Style = "Normal"

# Row 19
$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = "0.3427"
$ws.Range("E19").Value = "-0.75%"
$rng.Style = "Normal"

# Row 20
$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$ws.Range("E20").Value = "3.37%"
$rng.Style = "Normal"

# Row 21
$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "4.645"
$ws.Range("E21").Value = "-3.55%"
$rng.Style = "Normal"

# Row 22
$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$ws.Range("E22").Value = "-7.34%"
$rng.Style = "Normal"

# Row 23
$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.04382"
$ws.Range("E23").Value = "-2.21%"
$rng.Style = "Normal"

# Row 24
$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.001221"
$ws.Range("E24").Value = "-1.72%"
$rng.Style = "Normal"

# Row 25
$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.004338"
$ws.Range("E25").Value = "3.94%"
$rng.Style = "Normal"

# Row 27
$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = "0.0004005"
$ws.Range("E27").Value = "-9.96%"
$rng.Style = "Normal"

# Row 39
$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.02026"
$ws.Range("E39").Value = "-3.88%"
$rng.Style = "Normal"

# Row 40
$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.05043"
$ws.Range("E40").Value = "-2.41%"
$rng.Style = "Normal"

# Row 41
$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.007383"
$ws.Range("E41").Value = "-1.08%"
$rng.Style = "Normal"

# Row 42
$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.009750"
$ws.Range("E42").Value = "-3.75%"
$rng.Style = "Normal"

# Row 43
$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.1362"
$ws.Range("E43").Value = "-0.40%"
$rng.Style = "Normal"

# Row 44
$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.002133"
$ws.Range("E44").Value = "0.04%"
$rng.Style = "Normal"

# Row 45
$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.008715"
$ws.Range("E45").Value = "-10.02%"
$rng.Style = "Normal"

# Row 46
$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00006356"
$ws.Range("E46").Value = "0.44%"
$rng.Style = "Normal"

# Row 47
$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$ws.Range("E47").Value = "-0.03%"
$rng.Style = "Normal"

# Row 50
$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.03%"
$rng.Style = "Normal"

# Row 51
$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.03%"
$rng.Style = "Normal"
